$wb = $excel.ActiveWorkbook

# --- Sheet2 (validCredentialTest): move selection off the full-column selection, no longer the active tab ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A4").Select()

# --- Sheet3: rename "Sheet2" -> "addPatientTest" and populate the common dataprovider rows ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "addPatientTest"

# Write cell values in the exact order the strings were first introduced so the shared-string
# table layout matches the authored workbook.
$ws3.Range("A1").Value = "Username"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "Lanaguage"
$ws3.Range("D1").Value = "FirstName"
$ws3.Range("E1").Value = "LastName"
$ws3.Range("F1").Value = "DOB"
$ws3.Range("G1").Value = "Gender"
$ws3.Range("I1").Value = "ExpectedValue"

$ws3.Range("B2").Value = "pass"
$ws3.Range("C2").Value = "English (Indian)"
$ws3.Range("D2").Value = "John"
$ws3.Range("F2").Value = "2022-06-02"
$ws3.Range("A2").Value = "Admin"
$ws3.Range("G2").Value = "Male"
$ws3.Range("I2").Value = "John Wick"
$ws3.Range("E2").Value = "Wick"

$ws3.Range("H1").Value = "Expected Alert Text"
$ws3.Range("H2").Value = "Tobacco"

$ws3.Range("A3").Value = "accountant"
$ws3.Range("B3").Value = "accountant"
$ws3.Range("C3").Value = "English (Indian)"
$ws3.Range("D3").Value = "John"
$ws3.Range("E3").Value = "Wick"
$ws3.Range("F3").Value = "2022-06-02"
$ws3.Range("G3").Value = "Male"
$ws3.Range("H3").Value = "Tobacco"
$ws3.Range("I3").Value = "John Wick"

# Apply the "Text" (@) number format to the cells that carry it in the authored sheet
# (everything except C2, A3, B3, C3, which stay on the default/General style).
$ws3.Range("A1:I1,A2:B2,D2:I2,D3:I3").NumberFormat = "@"

# Column widths (best-fit sizing of the new dataprovider columns)
$ws3.Columns("A").ColumnWidth = 10.7109375
$ws3.Columns("B").ColumnWidth = 9.42578125
$ws3.Columns("C").ColumnWidth = 14.85546875
$ws3.Columns("D").ColumnWidth = 10.140625
$ws3.Columns("E").ColumnWidth = 9.7109375
$ws3.Columns("F").ColumnWidth = 10.42578125
$ws3.Columns("G").ColumnWidth = 7.5703125
$ws3.Columns("H").ColumnWidth = 19.140625
$ws3.Columns("I").ColumnWidth = 14.42578125

# Print orientation for the new sheet
$ws3.PageSetup.Orientation = 1

# Activate addPatientTest as the selected tab, with the reported selection.
$ws3.Activate()
$ws3.Range("B7").Select()
